# Update Receptor/Edge expression metrics with new TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 36.89194233333333
$ws.Range("N2").Value = 110.675827
$ws.Range("O2").Value = 0.3567095043190808
$ws.Range("P2").Value = 0.3567095043190809
$ws.Range("Q2").Value = 157.8289556604972
$ws.Range("R2").Value = 1420.460600944475
$ws.Range("S2").Value = 0.3446964233525793
$ws.Range("T2").Value = 0.3446964233525793

# Row 3
$ws.Range("M3").Value = 42.68037399999999
$ws.Range("O3").Value = 0.4126780562577495
$ws.Range("P3").Value = 0.4126780562577496
$ws.Range("Q3").Value = 182.5926863583166
$ws.Range("S3").Value = 0.3987800949113419
$ws.Range("T3").Value = 0.3987800949113419

# Row 4
$ws.Range("M4").Value = 23.85061433333334
$ws.Range("N4").Value = 71.55184300000001
$ws.Range("O4").Value = 0.2306124394231696
$ws.Range("P4").Value = 0.2306124394231696
$ws.Range("Q4").Value = 102.0363069550306
$ws.Range("R4").Value = 918.326762595275
$ws.Range("S4").Value = 0.222845991170098
$ws.Range("T4").Value = 0.222845991170098

# Row 5
$ws.Range("M5").Value = 36.89194233333333
$ws.Range("N5").Value = 110.675827
$ws.Range("O5").Value = 0.3567095043190808
$ws.Range("P5").Value = 0.3567095043190809
$ws.Range("Q5").Value = 5.500527115329445
$ws.Range("R5").Value = 49.504744037965
$ws.Range("S5").Value = 0.01201308096650157
$ws.Range("T5").Value = 0.01201308096650157

# Row 6
$ws.Range("M6").Value = 42.68037399999999
$ws.Range("O6").Value = 0.4126780562577495
$ws.Range("P6").Value = 0.4126780562577496
$ws.Range("R6").Value = 57.27215366498999
$ws.Range("S6").Value = 0.0138979613464077
$ws.Range("T6").Value = 0.0138979613464077

# Row 7
$ws.Range("M7").Value = 23.85061433333334
$ws.Range("N7").Value = 71.55184300000001
$ws.Range("O7").Value = 0.2306124394231696
$ws.Range("P7").Value = 0.2306124394231696
$ws.Range("R7").Value = 32.004781614685
$ws.Range("S7").Value = 0.00776644825307164
$ws.Range("T7").Value = 0.00776644825307164
